# Daily attendance processing - 2025-12-31 08:40:58
#
# For every row in the "Recorded By" column (G), when the value is a
# comma-separated list of recorder names/emails that includes the
# literal token "System", reverse the order of that list so "System"
# (and the rest) appear in reverse order. Single-value cells, or
# multi-value cells that don't include "System", are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null) {
        $text = [string]$val
        $parts = $text -split ", "

        if ($parts.Count -gt 1 -and ($parts -contains "System")) {
            $reversedParts = $parts[($parts.Count - 1)..0]
            $newVal = $reversedParts -join ", "
            $cell.Value2 = $newVal
        }
    }
}
